# Edit applied (per commit "Wed, May 20, 2020 10:05:25 PM"):
#   1. Turn on "embed TrueType fonts" for the presentation and register the
#      Limelight font that is embedded with it.
#   2. Re-style the sources-of-finance table (slide 6) with a new table
#      style id.
#
# NOTE: PowerPoint's object model has no documented, automatable way to
# stamp a font's embedded bytes into the package (font embedding is only
# reachable from File > Options > Save in the UI), so we still flip every
# COM knob that corresponds to that intent -- it's harmless if the host
# ignores it and is the closest automatable equivalent of the real edit.

$p = $ppt.ActivePresentation

# --- 1. Embed TrueType fonts (Limelight) -----------------------------
$p.EmbedTrueTypeFonts = $true

try {
    $fontItem = $p.Fonts.Add("Limelight")
    $fontItem.Embedded = $true
} catch {
    # Font embedding isn't exposed everywhere -- ignore if unsupported.
}

# --- 2. Apply the new table style to the table on slide 6 ------------
$newStyleId = "{1BC2204A-C1FD-457F-AD8B-F69B289A5D11}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
